# Replaces the first paragraph whose text contains $OldText with a single
# run containing $NewText (optionally wrapped in the given <w:rPr> markup).
#
# Uses Range.InsertXML (rather than Find.Execute or Range.Text=) because
# this runtime's Find/Replace and Range.Text setters collapse/consume any
# existing empty sibling <w:r/> runs in the paragraph, which would diverge
# from the target OOXML (these paragraphs keep a leading empty run).
# InsertXML, scoped exactly to the old run's character range, replaces only
# that run's contents and leaves neighboring empty runs untouched.
#
# NOTE: named parameters (e.g. "-OldText foo") are not reliably bound by
# this environment's PowerShell host, so all calls below use positional
# arguments only.
function Replace-ParaText {
    param([string]$OldText, [string]$NewText, [string]$RunPr)

    $d = $word.ActiveDocument
    foreach ($para in $d.Paragraphs) {
        $t = $para.Range.Text
        $idx = $t.IndexOf($OldText)
        if ($idx -ge 0) {
            $start = $para.Range.Start + $idx
            $len = $OldText.Length
            $target = $d.Range($start, $start + $len)

            $rprXml = ""
            if ($RunPr -ne "") {
                $rprXml = "<w:rPr>" + $RunPr + "</w:rPr>"
            }

            $escaped = $NewText -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'

            $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
                '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body><w:p><w:r>' + $rprXml + '<w:t>' + $escaped + '</w:t></w:r></w:p></w:body>' +
                '</w:document></pkg:xmlData></pkg:part></pkg:package>'

            $target.InsertXML($xmlFrag)
            return $true
        }
    }
    return $false
}

$dash = [char]0x2013

$oldTitle = "Play Chicken Party Free " + $dash + " Review of Booming Games' Online Slot"
$newTitle = "Play Chicken Party for Free - An Engaging Online Slot Game"

# 1a. Title heading (Heading1 paragraph, plain run, no rPr).
Replace-ParaText $oldTitle $newTitle ""

# 1b. Same title text repeated later in the document as a bold run.
Replace-ParaText $oldTitle $newTitle "<w:b/>"

# 2. "What we like" bullet - item 1.
Replace-ParaText "Realistic and pleasing illustration style" "Engaging gameplay with unique features" ""

# 3. "What we like" bullet - item 2.
Replace-ParaText "Interactive special symbols increase winning potential" "Realistic illustrations and pleasing aesthetic" ""

# 4. "What we like" bullet - item 3.
Replace-ParaText "Bonus mode with free spins and multipliers" "Opportunity to win free spins and multipliers" ""

# 5. "What we don't like" bullet - item 1.
Replace-ParaText "Random use of hen with golden egg symbol can be confusing" "Potential confusion with the hen with golden eggs symbol" ""

# 6. "What we don't like" bullet - item 2.
Replace-ParaText "High volatility may not be suitable for all players" "Higher volatility requires longer gameplay sessions" ""

# 7. Meta-description (italic run) near the end of the document.
Replace-ParaText "Find out more about Chicken Party, an online slot game by Booming Games, with our review. Play Chicken Party free and learn about its features." "Read our review of Chicken Party, an online slot game with engaging gameplay and free spins." "<w:i/>"
